$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.427.51"
$ws.Range("E2").Value = "  +3.96%  "
$ws.Range("D3").Value = "'1.837.53"
$ws.Range("E3").Value = "  +3.52%  "
$ws.Range("D4").Value = "'1.028"
$ws.Range("E4").Value = "  +2.58%  "
$ws.Range("D5").Value = "'318.29"
$ws.Range("E5").Value = "  +3.95%  "
$ws.Range("D6").Value = "'1.025"
$ws.Range("E6").Value = "  +2.26%  "
$ws.Range("D7").Value = "'0.4361"
$ws.Range("E7").Value = "  +3.02%  "
$ws.Range("D8").Value = "'0.3724"
$ws.Range("E8").Value = "  +3.32%  "
$ws.Range("D9").Value = "'0.07366"
$ws.Range("E9").Value = "  +3.42%  "
$ws.Range("D10").Value = "'0.8730"
$ws.Range("E10").Value = "  +4.47%  "
$ws.Range("D11").Value = "'21.39"
$ws.Range("E11").Value = "  +5.15%  "
$ws.Range("D12").Value = "'1.882.67"
$ws.Range("E12").Value = "  +5.84%  "
$ws.Range("D13").Value = "'5.461"
$ws.Range("E13").Value = "  +4.42%  "
$ws.Range("D14").Value = "'6.689"
$ws.Range("E14").Value = "  +3.86%  "
$ws.Range("D15").Value = "'0.07141"
$ws.Range("E15").Value = "  +3.99%  "
$ws.Range("D16").Value = "'82.61"
$ws.Range("E16").Value = "  +4.92%  "
$ws.Range("D17").Value = "'1.029"
$ws.Range("E17").Value = "  +2.29%  "
$ws.Range("D18").Value = "'0.000008984"
$ws.Range("E18").Value = "  +3.63%  "
$ws.Range("D19").Value = "'1.024"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").Value = "'15.40"
$ws.Range("E20").Value = "  +3.46%  "
$ws.Range("D21").Value = "'27.430.11"
$ws.Range("E21").Value = "  +3.93%  "
$ws.Range("D22").Value = "'5.235"
$ws.Range("E22").Value = "  +3.38%  "
$ws.Range("D23").Value = "'11.16"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("D24").Value = "'2.091.69"
$ws.Range("E24").Value = "  +4.76%  "
$ws.Range("D25").Value = "'156.76"
$ws.Range("E25").Value = "  +2.92%  "
$ws.Range("E26").Value = "  +6.05%  "
$ws.Range("D27").Value = "'18.59"
$ws.Range("E27").Value = "  +3.18%  "
$ws.Range("D28").Value = "'5.234"
$ws.Range("E28").Value = "  +3.42%  "
$ws.Range("D29").Value = "'1.923"
$ws.Range("E29").Value = "  +5.82%  "
$ws.Range("D30").Value = "'115.94"
$ws.Range("E30").Value = "  +1.76%  "
$ws.Range("D31").Value = "'0.09051"
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("D32").Value = "'1.201"
$ws.Range("E32").Value = "  +7.72%  "
$ws.Range("D33").Value = "'0.7606"
$ws.Range("E33").Value = "  +4.97%  "
$ws.Range("D34").Value = "'4.477"
$ws.Range("E34").Value = "  +3.89%  "
$ws.Range("D35").Value = "'2.867"
$ws.Range("E35").Value = "  +4.68%  "
$ws.Range("E36").Value = "  +2.43%  "
$ws.Range("D37").Value = "'1.145"
$ws.Range("E37").Value = "  +5.61%  "
$ws.Range("D38").Value = "'0.01964"
$ws.Range("E38").Value = "  +4.46%  "
$ws.Range("D39").Value = "'0.05247"
$ws.Range("E39").Value = "  +2.30%  "
$ws.Range("D40").Value = "'0.5167"
$ws.Range("E40").Value = "  +5.45%  "
$ws.Range("D41").Value = "'2.785"
$ws.Range("E41").Value = "  +6.88%  "
$ws.Range("D42").Value = "'0.1661"
$ws.Range("E42").Value = "  +3.31%  "
$ws.Range("D43").Value = "'6.556"
$ws.Range("E43").Value = "  +3.07%  "
$ws.Range("D44").Value = "'8.488"
$ws.Range("E44").Value = "  +6.73%  "
$ws.Range("D45").Value = "'108.91"
$ws.Range("E45").Value = "  +4.16%  "
$ws.Range("D46").Value = "'10.59"
$ws.Range("E46").Value = "  +3.92%  "
$ws.Range("D47").Value = "'1.027"
$ws.Range("E47").Value = "  +2.51%  "
$ws.Range("D48").Value = "'1.682"
$ws.Range("E48").Value = "  +2.89%  "
$ws.Range("D49").Value = "'0.4626"
$ws.Range("E49").Value = "  +4.49%  "
$ws.Range("D50").Value = "'1.896"
$ws.Range("E50").Value = "  +10.62%  "
$ws.Range("D51").Value = "'0.06306"
$ws.Range("E51").Value = "  +2.16%  "
